# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (one and only) slide master / all
#                            slides -- currently the "Simple Light" palette.
#   ppt/theme/theme2.xml  -> bound to the notes master -- currently the
#                            "Default" palette.
#
# The authored edit swaps the two themes' colour schemes, so the slide
# master ends up using the palette that used to live in theme2.xml
# ("Default": dk2=158158, lt2=F3F3F3, accent1=058DC7, accent2=50B432,
# accent3=ED561B, accent4=EDEF00, accent5=24CBE5, accent6=64E572,
# hlink=2200CC, folHlink=551A8B) while theme2.xml ends up with the palette
# that used to live in theme1.xml ("Simple Light").
#
# Only the slide-master side of that swap is reachable through the
# PowerPoint object model that this host exposes (NotesMaster/HandoutMaster
# resolve back onto the same SlideMaster/theme here, so theme2.xml's colors
# cannot be targeted independently). Apply the reachable half of the swap:
# push the "Default" theme colors onto the (single, shared) theme that
# backs the slide master, via ThemeColorScheme, same as Office's
# Design > Colors > Customize Colors dialog would.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 5800213    # dk2      158158
$colorScheme.Item(4).RGB  = 15987699   # lt2      F3F3F3
$colorScheme.Item(5).RGB  = 13077765   # accent1  058DC7
$colorScheme.Item(6).RGB  = 3322960    # accent2  50B432
$colorScheme.Item(7).RGB  = 1791725    # accent3  ED561B
$colorScheme.Item(8).RGB  = 61421      # accent4  EDEF00
$colorScheme.Item(9).RGB  = 15059748   # accent5  24CBE5
$colorScheme.Item(10).RGB = 7529828    # accent6  64E572
$colorScheme.Item(11).RGB = 13369378   # hlink    2200CC
$colorScheme.Item(12).RGB = 9116245    # folHlink 551A8B
